# Estadisticos Segundo Parcial 23 Mayo
#
# 1) "Estadisticos 2P" - fill in second-partial exam stats (previously a
#    copy of the totals-only placeholder) for both groups, plus the
#    Promedio (H) column which didn't exist yet.
# 2) "Estadisticos Final" - recompute the overall average (H) now that the
#    second partial has real numbers.
# 3) "Rescatables" - add two more students that need to re-take the exam
#    (ZUNO / TORRES families) ahead of the existing ANGEL ISMAEL VAZQUEZ
#    BONILLA record.

$wb = $excel.ActiveWorkbook

# --- Estadisticos 2P ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 4
$ws2.Range("F2").Value = 24
$ws2.Range("G2").Value = 85.70999999999999
$ws2.Range("H2").Value = 6.8

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 25
$ws2.Range("G3").Value = 100
$ws2.Range("H3").Value = 7.3

# --- Estadisticos Final --------------------------------------------------
$wsF = $wb.Worksheets.Item("Estadisticos Final")

$wsF.Range("H2").Value = 7.3
$wsF.Range("H3").Value = 7.8

# --- Rescatables -----------------------------------------------------
$wsR = $wb.Worksheets.Item("Rescatables")

# Two new students ahead of the existing one; the existing record (Angel
# Ismael Vazquez Bonilla) shifts from row 2 down to row 4.
$wsR.Range("A2").Value = 24330051920246
$wsR.Range("B2").Value = "ZUNO"
$wsR.Range("C2").Value = "FLORES"
$wsR.Range("D2").Value = "ALIN MARIEL"
$wsR.Range("E2").Value = "PREPARA SOLUCIONES PARA LAS OPERACIONES BÁSICAS DEL LABORATORIO"
$wsR.Range("F2").Value = "2ALCV"
$wsR.Range("G2").Value = 4

$wsR.Range("A3").Value = 24330051920238
$wsR.Range("B3").Value = "TORRES"
$wsR.Range("C3").Value = "PEREZ"
$wsR.Range("D3").Value = "ERIKA VALERIA"
$wsR.Range("E3").Value = "PREPARA SOLUCIONES PARA LAS OPERACIONES BÁSICAS DEL LABORATORIO"
$wsR.Range("F3").Value = "2ALCV"
$wsR.Range("G3").Value = 3

$wsR.Range("A4").Value = 24330051920239
$wsR.Range("B4").Value = "VAZQUEZ"
$wsR.Range("C4").Value = "BONILLA"
$wsR.Range("D4").Value = "ANGEL ISMAEL"
$wsR.Range("E4").Value = "PREPARA SOLUCIONES PARA LAS OPERACIONES BÁSICAS DEL LABORATORIO"
$wsR.Range("F4").Value = "2ALCV"
$wsR.Range("G4").Value = 3
